# predictionResults.xlsx — "Added multivariate series results for 16 or
# more customers": row 22 ("Min 11 purchases" bucket of the Customer Level
# Multivariate Series DayDiff2 table) previously had no evaluation data
# (B22/C22 were #DIV/0!, D22:K22 empty, L22 summed to 0). Fill in the newly
# available raw counts; the Precision/Recall ratios (B22, C22) and the
# Total (L22, a shared formula) recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D22").Value = 1041
$ws.Range("E22").Value = 2451
$ws.Range("F22").Value = 1061
$ws.Range("G22").Value = 2431
$ws.Range("H22").Value = 911
$ws.Range("I22").Value = 2301
$ws.Range("J22").Value = 981
$ws.Range("K22").Value = 3237

# Leave the cursor where the author left it after entering the data.
$null = $ws.Range("H23").Select()
